$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3559483367.572565
$ws.Range("C3").Value = 3565065203.8173027
$ws.Range("C4").Value = 3576228911.510071
$ws.Range("C5").Value = 3589625422.6994276
